# Update the "Förändrad" (Changed) date column (C) for every data row
# from 45179 (2023-09-10) to 45180 (2023-09-11), matching the automatic
# update recorded in the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 295

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45179) {
        $cell.Value = 45180
    }
}
